# Rebuild the document body paragraph from a flat list of
# (text, color) runs, as if driven by lines read from a plain-text
# input file (one colored word/token at a time), per the "Using a txt
# as input" commit.
#
# Each element is a run: Text is the literal run text (already
# including its trailing separator, e.g. a trailing space), Color is
# the RRGGBB hex string for w:color, $null for "no explicit color"
# (automatic), or the sentinel "BREAK" for a manual line break
# (<w:br/>) inside the same paragraph.

function HexToWordColor($hex) {
    $rr = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $gg = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $bb = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $rr + ($gg * 256) + ($bb * 65536)
}

$runs = @(
    @{ Text = "I "; Color = "7FCDE7" }
    @{ Text = "can "; Color = "999F54" }
    @{ Text = "not "; Color = "6C913F" }
    @{ Text = "stomach "; Color = "3EAFF3" }
    @{ Text = "these "; Color = "D39CC2" }
    @{ Text = "forms "; Color = "BB50F1" }
    @{ Text = "and "; Color = "D9E22C" }
    @{ Text = "colors "; Color = "BB50F1" }
    @{ Text = "anymore "; Color = "F1AAAE" }
    @{ Text = ". "; Color = "1946EA" }
    @{ Text = "- "; Color = $null }
    @{ Text = "Ya "; Color = "F1AAAE" }
    @{ Text = "no "; Color = "F1AAAE" }
    @{ Text = "puedo "; Color = "999F54" }
    @{ Text = "soportar "; Color = "3EAFF3" }
    @{ Text = "estas "; Color = "D39CC2" }
    @{ Text = "formas "; Color = "BB50F1" }
    @{ Text = "y "; Color = "D9E22C" }
    @{ Text = "colores "; Color = "BB50F1" }
    @{ Text = ". "; Color = "1946EA" }
    @{ Text = ""; Color = "BREAK" }
    @{ Text = "But "; Color = "D9E22C" }
    @{ Text = "I "; Color = "7FCDE7" }
    @{ Text = "'m "; Color = "999F54" }
    @{ Text = "here "; Color = "F1AAAE" }
    @{ Text = "to "; Color = "6C913F" }
    @{ Text = "continue "; Color = "3EAFF3" }
    @{ Text = "after "; Color = "F1AAAE" }
    @{ Text = "all "; Color = "F1AAAE" }
    @{ Text = "I "; Color = "7FCDE7" }
    @{ Text = "have "; Color = "999F54" }
    @{ Text = "been "; Color = "999F54" }
    @{ Text = "through "; Color = "F4EAF4" }
    @{ Text = ". "; Color = "1946EA" }
    @{ Text = "- "; Color = $null }
    @{ Text = "Pero "; Color = "D9E22C" }
    @{ Text = "estoy "; Color = "3EAFF3" }
    @{ Text = "aquí "; Color = "F1AAAE" }
    @{ Text = "para "; Color = "F4EAF4" }
    @{ Text = "continuar "; Color = "3EAFF3" }
    @{ Text = "después "; Color = "F1AAAE" }
    @{ Text = "de "; Color = "F4EAF4" }
    @{ Text = "todo "; Color = "D39CC2" }
    @{ Text = "lo "; Color = "7FCDE7" }
    @{ Text = "que "; Color = "7FCDE7" }
    @{ Text = "he "; Color = "999F54" }
    @{ Text = "pasado "; Color = "3EAFF3" }
    @{ Text = ". "; Color = "1946EA" }
    @{ Text = ""; Color = "BREAK" }
    @{ Text = "I "; Color = "7FCDE7" }
    @{ Text = "try "; Color = "3EAFF3" }
    @{ Text = "to "; Color = "6C913F" }
    @{ Text = "keep "; Color = "3EAFF3" }
    @{ Text = "my "; Color = "7FCDE7" }
    @{ Text = "eyes "; Color = "BB50F1" }
    @{ Text = "open "; Color = "55F12C" }
    @{ Text = ", "; Color = "1946EA" }
    @{ Text = "I "; Color = "7FCDE7" }
    @{ Text = "am "; Color = "999F54" }
    @{ Text = "realizing "; Color = "3EAFF3" }
    @{ Text = ". "; Color = "1946EA" }
    @{ Text = "- "; Color = $null }
    @{ Text = "Intento "; Color = "3EAFF3" }
    @{ Text = "mantener "; Color = "3EAFF3" }
    @{ Text = "los "; Color = "D39CC2" }
    @{ Text = "ojos "; Color = "BB50F1" }
    @{ Text = "abiertos "; Color = "55F12C" }
    @{ Text = ", "; Color = "1946EA" }
    @{ Text = "me "; Color = "7FCDE7" }
    @{ Text = "estoy "; Color = "999F54" }
    @{ Text = "dando "; Color = "3EAFF3" }
    @{ Text = "cuenta "; Color = "BB50F1" }
    @{ Text = ". "; Color = "1946EA" }
    @{ Text = ""; Color = "BREAK" }
    @{ Text = "This "; Color = "D39CC2" }
    @{ Text = "life "; Color = "BB50F1" }
    @{ Text = "and "; Color = "D9E22C" }
    @{ Text = "death "; Color = "BB50F1" }
    @{ Text = "more "; Color = "F1AAAE" }
    @{ Text = "precious "; Color = "55F12C" }
    @{ Text = "than "; Color = "F4EAF4" }
    @{ Text = "anything "; Color = "7FCDE7" }
    @{ Text = "- "; Color = $null }
    @{ Text = "Esta "; Color = "D39CC2" }
    @{ Text = "vida "; Color = "BB50F1" }
    @{ Text = "y "; Color = "D9E22C" }
    @{ Text = "muerte "; Color = "BB50F1" }
    @{ Text = "más "; Color = "F1AAAE" }
    @{ Text = "preciosa "; Color = "55F12C" }
    @{ Text = "que "; Color = "007F67" }
    @{ Text = "cualquier "; Color = "D39CC2" }
    @{ Text = "cosa "; Color = "BB50F1" }
    @{ Text = ""; Color = "BREAK" }
)

$d = $word.ActiveDocument

# Wipe the existing paragraph content; we're rebuilding it from scratch.
$d.Content.Delete()

$pos = 0
foreach ($run in $runs) {
    $ip = $d.Range($pos, $pos)

    if ($run.Color -eq "BREAK") {
        $ip.InsertBreak(6)   # wdLineBreak
        $pos = $pos + 1
        continue
    }

    $ip.InsertAfter($run.Text)
    $len = $run.Text.Length

    if ($run.Color -ne $null) {
        $colored = $d.Range($pos, $pos + $len)
        $colored.Font.Color = HexToWordColor($run.Color)
    }

    $pos = $pos + $len
}
